$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34: num_customers 90 -> 91, retention_rate recalculated (91/2256)
$ws.Range("C34").Value = 91
$ws.Range("E34").Value = 91/2256

# Row 36: num_customers 148 -> 149, retention_rate recalculated (149/1930)
$ws.Range("C36").Value = 149
$ws.Range("E36").Value = 149/1930

# Row 37: num_customers and cohort_size 1003 -> 1014
$ws.Range("C37").Value = 1014
$ws.Range("D37").Value = 1014
